# US-14491 [IMP] PO Catalogue tab: Added % deviation at header and line level,
# price is bold if higher than the catalogue, fixed a column name, and added
# deviation to the mismatch report.
#
# This script inserts a new "% Price Deviation" column right before the
# existing "Catalogue SoQ" column (column M) on the PO Catalogue Mismatch
# sheet, shifting every column from M onward one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column M (13th column, "Catalogue SoQ").
# Excel copies the formatting of the column to the left, which matches the
# styling already used by the rest of the header/body cells in this report.
$ws.Columns("M:M").Insert()

# Populate the new header cell with the new label.
$ws.Range("M10").Value = "% Price Deviation"

# Match the new active selection recorded for the sheet view.
$ws.Range("M11").Select() | Out-Null
